$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 5 - "No adquirir la capasitacion..." (android training) risk
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = "No adquirir la capasitación necesaria para el desarrollo de aplicaciónes móviles "
$ws.Range("E5").Value = "Desfase de tiempo para terminar el proyecto(la parte de moviles)"
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = "Aprender lo mas amplio posible del manejo de android en horas extraclase"
$ws.Range("J5").Value = "Se tendrá que posponer la fecha de entrega."
$ws.Rows.Item(5).RowHeight = 25.5

# ---------------------------------------------------------------------------
# Row 6 - web graphics risk
# ---------------------------------------------------------------------------
$ws.Range("D6").Value = "No saber manipular correctamente los gráficos en la web."
$ws.Range("E6").Value = "No terminar el proyecto completo(la parte de web)"
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = "Aprender a desarrollar aplicaciones web en horas extraclase."
$ws.Range("J6").Value = "No se podrá terminar ningún juego web."
$ws.Rows.Item(6).RowHeight = 12.75

# ---------------------------------------------------------------------------
# Row 7 - dependent tasks risk
# ---------------------------------------------------------------------------
$ws.Range("D7").Value = "Alguna de las tareas que tienen dependencia no sea terminada en el tiempo planeado."
$ws.Range("E7").Value = "Se retrasarán los avances de nuevas tareas."
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = "Dar prioridad a las tareas con dependencia."
$ws.Range("J7").Value = "Finalizar la tareas lo más pronto posible, para continuar con la tarea que se debe de realizar después."

# ---------------------------------------------------------------------------
# Row 8 - new risk: client rejects prototypes
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "Nuestro cliente no acepte los prototipos"
$ws.Range("E8").Value = "Desfase de tiempo para rediseñarlos."
$ws.Range("G8").Value = "Baja"
$ws.Range("I8").Value = "Hacer que los prototipos cumplan con la mayoría de los requisitos propuestos por el cliente."
$ws.Range("J8").Value = "Volver a hacer nuevos prototipos."
$ws.Rows.Item(8).RowHeight = 25.5

# D8/E8 lose the left/center alignment that the rest of the column uses,
# keeping only wrap text (border stays the thin box it already had).
$ws.Range("D8:E8").HorizontalAlignment = 1
$ws.Range("D8:E8").VerticalAlignment = -4107
$ws.Range("D8:E8").WrapText = $true

# F8/H8 pick up the same full thin-box border that G8 already has.
$ws.Range("F8").Borders.Item(7).LineStyle = 1
$ws.Range("F8").Borders.Item(8).LineStyle = 1
$ws.Range("F8").Borders.Item(9).LineStyle = 1
$ws.Range("F8").Borders.Item(10).LineStyle = 1
$ws.Range("H8").Borders.Item(7).LineStyle = 1
$ws.Range("H8").Borders.Item(8).LineStyle = 1
$ws.Range("H8").Borders.Item(9).LineStyle = 1
$ws.Range("H8").Borders.Item(10).LineStyle = 1

# I8 drops its left border and left/center alignment, keeping only wrap.
$ws.Range("I8").Borders.Item(7).LineStyle = -4142
$ws.Range("I8").HorizontalAlignment = 1
$ws.Range("I8").VerticalAlignment = -4107
$ws.Range("I8").WrapText = $true

# ---------------------------------------------------------------------------
# Selection moves to E8 (last cell edited by the author)
# ---------------------------------------------------------------------------
$ws.Range("E8").Select()
